$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Result) and E (Entered Date) would otherwise be auto-converted
# by Excel's type inference (e.g. "100.00" -> 100, "2024-10-09" -> a date
# serial). Mark them as Text first so the literal strings are preserved,
# then restore the default "Normal" style so no custom number format is
# left applied to the new cells.
$textRange = $ws.Range("D44:E45")
$textRange.NumberFormat = "@"

$ws.Range("A44").Value = "2024-10-09 22:32:03"
$ws.Range("B44").Value = "get_price"
$ws.Range("C44").Value = "https://example.com/product"
$ws.Range("D44").Value = "100.00"
$ws.Range("E44").Value = "2024-10-09"
$ws.Range("F44").Value = "22:32:03"

$ws.Range("A45").Value = "2024-10-09 23:31:35"
$ws.Range("B45").Value = "get_price"
$ws.Range("C45").Value = "https://example.com/product"
$ws.Range("D45").Value = "100.00"
$ws.Range("E45").Value = "2024-10-09"
$ws.Range("F45").Value = "23:31:35"

$textRange.Style = "Normal"
